# Update TPM-derived metrics for Hbegf-Cd82 LR pair (YoungD7) with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.082188
$ws.Range("H2").Value = 24.246564
$ws.Range("I2").Value = 0.2755301789948819
$ws.Range("J2").Value = 0.2755301789948819
$ws.Range("M2").Value = 3.932154333333333
$ws.Range("N2").Value = 11.796463
$ws.Range("O2").Value = 0.04140655843753142
$ws.Range("P2").Value = 0.04140655843753142
$ws.Range("Q2").Value = 31.78041056701467
$ws.Range("R2").Value = 286.023695103132
$ws.Range("S2").Value = 0.01140875645785507
$ws.Range("T2").Value = 0.01140875645785507
$ws.Range("G3").Value = 8.082188
$ws.Range("H3").Value = 24.246564
$ws.Range("I3").Value = 0.2755301789948819
$ws.Range("J3").Value = 0.2755301789948819
$ws.Range("O3").Value = 0.008782731887949688
$ws.Range("P3").Value = 0.008782731887949688
$ws.Range("Q3").Value = 6.740932737024001
$ws.Range("R3").Value = 60.66839463321602
$ws.Range("S3").Value = 0.002419907689150834
$ws.Range("T3").Value = 0.002419907689150835
$ws.Range("G4").Value = 8.082188
$ws.Range("H4").Value = 24.246564
$ws.Range("I4").Value = 0.2755301789948819
$ws.Range("J4").Value = 0.2755301789948819
$ws.Range("M4").Value = 79.04521166666666
$ws.Range("N4").Value = 237.135635
$ws.Range("O4").Value = 0.8323656445367243
$ws.Range("P4").Value = 0.8323656445367243
$ws.Range("Q4").Value = 638.8582611897933
$ws.Range("R4").Value = 5749.72435070814
$ws.Range("S4").Value = 0.2293418550283939
$ws.Range("T4").Value = 0.2293418550283939
$ws.Range("G5").Value = 8.082188
$ws.Range("H5").Value = 24.246564
$ws.Range("I5").Value = 0.2755301789948819
$ws.Range("J5").Value = 0.2755301789948819
$ws.Range("M5").Value = 11.15311533333333
$ws.Range("N5").Value = 33.459346
$ws.Range("O5").Value = 0.1174450651377945
$ws.Range("P5").Value = 0.1174450651377945
$ws.Range("Q5").Value = 90.14157490968266
$ws.Range("R5").Value = 811.274174187144
$ws.Range("S5").Value = 0.03235965981948209
$ws.Range("T5").Value = 0.03235965981948209
$ws.Range("I6").Value = 0.4533445791334642
$ws.Range("J6").Value = 0.4533445791334642
$ws.Range("M6").Value = 3.932154333333333
$ws.Range("N6").Value = 11.796463
$ws.Range("O6").Value = 0.04140655843753142
$ws.Range("P6").Value = 0.04140655843753142
$ws.Range("Q6").Value = 52.29001376818177
$ws.Range("R6").Value = 470.610123913636
$ws.Range("S6").Value = 0.01877143880822787
$ws.Range("T6").Value = 0.01877143880822787
$ws.Range("I7").Value = 0.4533445791334642
$ws.Range("J7").Value = 0.4533445791334642
$ws.Range("O7").Value = 0.008782731887949688
$ws.Range("P7").Value = 0.008782731887949688
$ws.Range("S7").Value = 0.003981603891384606
$ws.Range("T7").Value = 0.003981603891384606
$ws.Range("I8").Value = 0.4533445791334642
$ws.Range("J8").Value = 0.4533445791334642
$ws.Range("M8").Value = 79.04521166666666
$ws.Range("N8").Value = 237.135635
$ws.Range("O8").Value = 0.8323656445367243
$ws.Range("P8").Value = 0.8323656445367243
$ws.Range("Q8").Value = 1051.147756668802
$ws.Range("R8").Value = 9460.329810019219
$ws.Range("S8").Value = 0.3773484528076559
$ws.Range("T8").Value = 0.3773484528076559
$ws.Range("I9").Value = 0.4533445791334642
$ws.Range("J9").Value = 0.4533445791334642
$ws.Range("M9").Value = 11.15311533333333
$ws.Range("N9").Value = 33.459346
$ws.Range("O9").Value = 0.1174450651377945
$ws.Range("P9").Value = 0.1174450651377945
$ws.Range("Q9").Value = 148.3147671479458
$ws.Range("R9").Value = 1334.832904331512
$ws.Range("S9").Value = 0.05324308362619574
$ws.Range("T9").Value = 0.05324308362619575
$ws.Range("G10").Value = 5.789497666666667
$ws.Range("H10").Value = 17.368493
$ws.Range("I10").Value = 0.1973699855023315
$ws.Range("J10").Value = 0.1973699855023315
$ws.Range("M10").Value = 3.932154333333333
$ws.Range("N10").Value = 11.796463
$ws.Range("O10").Value = 0.04140655843753142
$ws.Range("P10").Value = 0.04140655843753142
$ws.Range("Q10").Value = 22.76519833780655
$ws.Range("R10").Value = 204.886785040259
$ws.Range("S10").Value = 0.008172411838517017
$ws.Range("T10").Value = 0.008172411838517019
$ws.Range("G11").Value = 5.789497666666667
$ws.Range("H11").Value = 17.368493
$ws.Range("I11").Value = 0.1973699855023315
$ws.Range("J11").Value = 0.1973699855023315
$ws.Range("O11").Value = 0.008782731887949688
$ws.Range("P11").Value = 0.008782731887949688
$ws.Range("Q11").Value = 4.828718949888001
$ws.Range("R11").Value = 43.45847054899201
$ws.Range("S11").Value = 0.001733447665395494
$ws.Range("T11").Value = 0.001733447665395495
$ws.Range("G12").Value = 5.789497666666667
$ws.Range("H12").Value = 17.368493
$ws.Range("I12").Value = 0.1973699855023315
$ws.Range("J12").Value = 0.1973699855023315
$ws.Range("M12").Value = 79.04521166666666
$ws.Range("N12").Value = 237.135635
$ws.Range("O12").Value = 0.8323656445367243
$ws.Range("P12").Value = 0.8323656445367243
$ws.Range("Q12").Value = 457.6320685053394
$ws.Range("R12").Value = 4118.688616548055
$ws.Range("S12").Value = 0.1642839951948521
$ws.Range("T12").Value = 0.1642839951948521
$ws.Range("G13").Value = 5.789497666666667
$ws.Range("H13").Value = 17.368493
$ws.Range("I13").Value = 0.1973699855023315
$ws.Range("J13").Value = 0.1973699855023315
$ws.Range("M13").Value = 11.15311533333333
$ws.Range("N13").Value = 33.459346
$ws.Range("O13").Value = 0.1174450651377945
$ws.Range("P13").Value = 0.1174450651377945
$ws.Range("Q13").Value = 64.57093519839755
$ws.Range("R13").Value = 581.138416785578
$ws.Range("S13").Value = 0.02318013080356688
$ws.Range("T13").Value = 0.02318013080356689
$ws.Range("G14").Value = 2.163479333333334
$ws.Range("H14").Value = 6.490438
$ws.Range("I14").Value = 0.0737552563693224
$ws.Range("J14").Value = 0.0737552563693224
$ws.Range("M14").Value = 3.932154333333333
$ws.Range("N14").Value = 11.796463
$ws.Range("O14").Value = 0.04140655843753142
$ws.Range("P14").Value = 0.04140655843753142
$ws.Range("Q14").Value = 8.507134635643778
$ws.Range("R14").Value = 76.56421172079399
$ws.Range("S14").Value = 0.003053951332931459
$ws.Range("T14").Value = 0.00305395133293146
$ws.Range("G15").Value = 2.163479333333334
$ws.Range("H15").Value = 6.490438
$ws.Range("I15").Value = 0.0737552563693224
$ws.Range("J15").Value = 0.0737552563693224
$ws.Range("O15").Value = 0.008782731887949688
$ws.Range("P15").Value = 0.008782731887949688
$ws.Range("Q15").Value = 1.804445611008
$ws.Range("R15").Value = 16.240010499072
$ws.Range("S15").Value = 0.0006477726420187522
$ws.Range("T15").Value = 0.0006477726420187522
$ws.Range("G16").Value = 2.163479333333334
$ws.Range("H16").Value = 6.490438
$ws.Range("I16").Value = 0.0737552563693224
$ws.Range("J16").Value = 0.0737552563693224
$ws.Range("M16").Value = 79.04521166666666
$ws.Range("N16").Value = 237.135635
$ws.Range("O16").Value = 0.8323656445367243
$ws.Range("P16").Value = 0.8323656445367243
$ws.Range("Q16").Value = 171.0126818397922
$ws.Range("R16").Value = 1539.11413655813
$ws.Range("S16").Value = 0.06139134150582238
$ws.Range("T16").Value = 0.06139134150582238
$ws.Range("G17").Value = 2.163479333333334
$ws.Range("H17").Value = 6.490438
$ws.Range("I17").Value = 0.0737552563693224
$ws.Range("J17").Value = 0.0737552563693224
$ws.Range("M17").Value = 11.15311533333333
$ws.Range("N17").Value = 33.459346
$ws.Range("O17").Value = 0.1174450651377945
$ws.Range("P17").Value = 0.1174450651377945
$ws.Range("Q17").Value = 24.12953452594978
$ws.Range("R17").Value = 217.165810733548
$ws.Range("S17").Value = 0.008662190888549803
$ws.Range("T17").Value = 0.008662190888549805

$wb.Save()
